# Mise à jour de l'application
# Adds 9 new training-log rows (533-541, all dated 2025-11-04 / serial 45965)
# to the bottom of the "Feuil1" sheet, extending the data table from
# A1:I532 to A1:I541, and recomputes the "Charge" (I) column = Volume*Intensite.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$lastRow = 532
$firstNew = 533
$lastNew = 541

# 1) Copy the formatting of the last existing row down onto the new rows so
#    that number formats / styles (date column, grey "Helvetica Neue" font,
#    centered empty "Localisation douleur" cells, etc.) match the rest of
#    the table.
$ws.Range("A$lastRow`:I$lastRow").Copy()
$ws.Range("A$firstNew`:I$lastNew").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 2) New row data, in the same column order as the sheet:
#    Date | Nom du joueur | Volume | Intensite | Fatigue | Douleur | Localisation douleur | Plaisir
$rows = @(
    @{ Row=533; Date=45965; Nom="Levy Ndoutoume";     Volume=70; Intensite=7; Fatigue=8; Douleur=1; Loc="Ischio";     Plaisir=5 },
    @{ Row=534; Date=45965; Nom="Yoann Martelat";      Volume=70; Intensite=7; Fatigue=6; Douleur=6; Loc="Genou";      Plaisir=7 },
    @{ Row=535; Date=45965; Nom="Maé Clavel";          Volume=70; Intensite=6; Fatigue=4; Douleur=6; Loc="Ischio";     Plaisir=7 },
    @{ Row=536; Date=45965; Nom="Naim Ighbane";        Volume=70; Intensite=7; Fatigue=5; Douleur=0; Loc=$null;       Plaisir=4 },
    @{ Row=537; Date=45965; Nom="Karim Belmahi";       Volume=70; Intensite=5; Fatigue=7; Douleur=0; Loc=$null;       Plaisir=10 },
    @{ Row=538; Date=45965; Nom="Ilan Ihaddadene";     Volume=70; Intensite=8; Fatigue=7; Douleur=0; Loc=$null;       Plaisir=9 },
    @{ Row=539; Date=45965; Nom="Karahali Souaré";     Volume=70; Intensite=6; Fatigue=7; Douleur=7; Loc="Cheville";   Plaisir=4 },
    @{ Row=540; Date=45965; Nom="Naim Dhib";           Volume=70; Intensite=6; Fatigue=7; Douleur=1; Loc="Courbature"; Plaisir=7 },
    @{ Row=541; Date=45965; Nom="Sofiane Belle";       Volume=70; Intensite=7; Fatigue=7; Douleur=7; Loc="Ischio";     Plaisir=5 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.Nom
    $ws.Cells.Item($row, 3).Value = $r.Volume
    $ws.Cells.Item($row, 4).Value = $r.Intensite
    $ws.Cells.Item($row, 5).Value = $r.Fatigue
    $ws.Cells.Item($row, 6).Value = $r.Douleur
    if ($r.Loc) {
        $ws.Cells.Item($row, 7).Value = $r.Loc
    } else {
        $ws.Cells.Item($row, 7).ClearContents()
    }
    $ws.Cells.Item($row, 8).Value = $r.Plaisir
    $ws.Cells.Item($row, 9).Formula = "=C$row*D$row"
}

# 3) The "Localisation douleur" column for the three rows with no value
#    should look exactly like the other empty cells in that column
#    (centered style, no content) - re-stamp their format to be sure.
$ws.Range("G2").Copy()
$ws.Range("G536").PasteSpecial(-4122)
$ws.Range("G537").PasteSpecial(-4122)
$ws.Range("G538").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$excel.Calculate()

# 4) Update the view so the selection / scroll position follow the newly
#    added rows, like Excel would after entering this data.
$ws.Range("A$firstNew").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 514
$win.ScrollColumn = 1
$ws.Range("L535").Select()

Write-Host "Added rows $firstNew to $lastNew; dimension is now $($ws.UsedRange.Address())"
